$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Overview" sheet - b.md row (row 3) moved from "Handed back" to
# "Ready for handoff", with an updated handoff generation timestamp.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-04 00:41:58"

# ---------------------------------------------------------------------------
# "zh-cn" sheet - b.md row (row 3): new handback xliff generated, status
# flips to "Ready for handoff", duplicate flag clears, and an error detail
# about the handback not matching the latest source version is recorded.
#
# NOTE: the literal text "False" would normally be auto-coerced to a real
# boolean by plain assignment (mirrors Excel's own "typed input" parsing).
# The source workbook stores it as literal text instead, so we force text
# with a leading apostrophe and then reset the style back to Normal to
# drop the quote-prefix marker that the apostrophe trick leaves behind.
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").Style = "Normal"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-09-04 00:41:53"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8732e4d88870f38f60d210e592bc8b9fe17462ab/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9839233e08e228fc912ac089251c6063f29e8c68/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# "de-de" sheet - b.md row (row 3): same kind of update as zh-cn above.
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("F3").Style = "Normal"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-09-04 00:41:58"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8732e4d88870f38f60d210e592bc8b9fe17462ab/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9839233e08e228fc912ac089251c6063f29e8c68/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 39.17
